# Generate Report for Handback
# Updates the Overview status text, and fills in the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns (with
# hyperlinks) on the per-locale (zh-cn / de-de) sheets, now that the
# localized content has been handed back and is in sync with en-US.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: status goes from "Ready for handoff" to "Handed back:
# in sync with en-US" for both locales.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 5).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(2, 6).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.166666666666664
$overview.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column mirrors the Overview sheet's wording.
$zhcn.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$zhcn.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

$zhcn.Cells.Item(2, 9).Value = "0a644613-e7b7-4934-ae13-de75f45037d7.md"
$zhcn.Cells.Item(2, 10).Value = "0a644613-e7b7-4934-ae13-de75f45037d7.c64eccfa866dee4d13e20d4dc6c9018b72dcbf13.zh-cn.xlf"
$zhcn.Cells.Item(3, 9).Value = "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md"
$zhcn.Cells.Item(3, 10).Value = "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.670553798e97d3d322576a330320f1f12ff9e194.zh-cn.xlf"

# Latest Handback DateTime (same value already shown by K2/K3 via the
# shared "0001-01-01 00:00:00" placeholder -- now the real handback time).
$zhcn.Cells.Item(2, 11).Value = "2016-09-06 03:09:52"
$zhcn.Cells.Item(3, 11).Value = "2016-09-06 03:09:52"

$zhcn.Hyperlinks.Add($zhcn.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8efe5b823709155fdd186ad3236d1f801c2bdb4f/e2e/0a644613-e7b7-4934-ae13-de75f45037d7.md", "", "", "0a644613-e7b7-4934-ae13-de75f45037d7.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8efe5b823709155fdd186ad3236d1f801c2bdb4f/e2e/ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md", "", "", "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md") | Out-Null

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666664
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column mirrors the Overview sheet's wording.
$dede.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$dede.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

$dede.Cells.Item(2, 9).Value = "0a644613-e7b7-4934-ae13-de75f45037d7.md"
$dede.Cells.Item(2, 10).Value = "0a644613-e7b7-4934-ae13-de75f45037d7.c64eccfa866dee4d13e20d4dc6c9018b72dcbf13.de-de.xlf"
$dede.Cells.Item(3, 9).Value = "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md"
$dede.Cells.Item(3, 10).Value = "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.670553798e97d3d322576a330320f1f12ff9e194.de-de.xlf"

# Latest Handback DateTime - de-de xliffs just came back, a few minutes
# after the zh-cn ones.
$dede.Cells.Item(2, 11).Value = "2016-09-06 03:10:00"
$dede.Cells.Item(3, 11).Value = "2016-09-06 03:10:00"

$dede.Hyperlinks.Add($dede.Cells.Item(2, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8efe5b823709155fdd186ad3236d1f801c2bdb4f/e2e/0a644613-e7b7-4934-ae13-de75f45037d7.md", "", "", "0a644613-e7b7-4934-ae13-de75f45037d7.md") | Out-Null
$dede.Hyperlinks.Add($dede.Cells.Item(3, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8efe5b823709155fdd186ad3236d1f801c2bdb4f/e2e/ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md", "", "", "ca695cbc-cc4f-4b5c-98d7-7f3918e272bd.md") | Out-Null

$dede.Columns.Item(3).ColumnWidth = 29.166666666666664
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
